# Update Name of Algo
# Applies the cell-value corrections shown in the diff for Sheet1
# (columns B and D of the result_data_RandomForest data).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("D4").Value = -7.862900000000004
$ws.Range("B9").Value = 5.553000000000005
$ws.Range("D9").Value = -7.502399999999995
$ws.Range("D11").Value = -7.887599999999997
$ws.Range("B18").Value = 6.485399999999998
$ws.Range("B20").Value = 8.9139
$ws.Range("D23").Value = -8.028100000000002
$ws.Range("D24").Value = -7.511499999999999
$ws.Range("D26").Value = -7.4418
$ws.Range("B27").Value = 6.535800000000005
$ws.Range("D34").Value = -7.786800000000004
$ws.Range("B35").Value = 8.495900000000008
$ws.Range("D35").Value = -8.016100000000003
$ws.Range("D48").Value = -7.505199999999999
$ws.Range("D49").Value = -8.086900000000005
$ws.Range("D52").Value = -7.923900000000003
$ws.Range("D66").Value = -7.221599999999999
$ws.Range("D67").Value = -6.876399999999997
$ws.Range("B69").Value = 5.387599999999992
$ws.Range("B76").Value = 5.503399999999998
$ws.Range("B78").Value = 8.625000000000007
$ws.Range("D78").Value = -7.843999999999999
$ws.Range("D80").Value = -7.8485
$ws.Range("B82").Value = 5.445400000000001
$ws.Range("B83").Value = 5.177799999999996
$ws.Range("B93").Value = 5.470999999999998
$ws.Range("D99").Value = -8.006700000000002
$ws.Range("D104").Value = -7.479
